$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells stay text-typed (matches source inlineStr cells) so numeric-looking
# strings such as "209.50" or "26.195.20" are not coerced into floating point numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.195.20"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.581.66"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.39%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.50"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.35%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.44%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.47"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.803.75"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.23%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.622.28"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.06"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.49"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.193.59"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.89%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.25"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "207.20"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.45%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.88"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.80"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.52%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.97%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.112"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.17%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.24%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.14"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.22"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.29%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.279.33"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.47"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.608"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.96%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.48"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.16"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.60%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.814"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.79%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.30%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.764"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.17%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.25"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.718.00"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.88"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.97%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.84%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.95%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.42"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.27%  "
